# Updated 90Zr reaction rates to account for the error in atom density
#
# 1) Row 7 was a stray duplicate entry (label "ENDF 115In(n,g)") that
#    needs to be removed entirely. Deleting the whole row shifts every
#    row below it up by one (rows 8-20 -> 7-19, the blank spacer row
#    21 -> 20, and so on down to the bottom of the sheet), and Excel
#    automatically drops the now-unused "ENDF 115In(n,g)" shared string.
#
# 2) The F2 cross-section value needs to be corrected for the atom
#    density error: multiply the original value by the ratio of the
#    correct to erroneous atom densities (0.022148/0.02551). This is
#    entered as a formula so the dependent P2 cell recalculates too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 7 entirely (shifts rows 8+ up by one row).
$ws.Rows(7).Delete()

# Correct the F2 cross-section for the atom-density error.
$ws.Range("F2").Formula = "=0.00000012365*0.022148/0.02551"

# Match the saved selection state of the edited workbook.
$ws.Range("A6").Select() | Out-Null
